# Update emission activity ratios for electricity generation technologies:
# add a new "EmissionActivityRatio" column (E) to the TFEC sheet and
# populate it for the relevant technology rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell E1, formatted like the neighbouring header cell D1
# (bold, no fill/border - same look as the other header cells).
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("E1").Value = "EmissionActivityRatio"

# EmissionActivityRatio values for the rows that have them
$ws.Range("E2").Value = 0.0961
$ws.Range("E3").Value = 0.0561
$ws.Range("E4").Value = 0.0741
$ws.Range("E5").Value = 0.1225
$ws.Range("E10").Value = 0.0741
$ws.Range("E11").Value = 0.0741

# Scroll the sheet so column C is the first visible column, and leave the
# selection on the last edited cell (E11), matching the author's view state.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E11").Select()
